$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2-230.
# The serial value 45182 (2023-09-13) was updated to 45184 (2023-09-15).
$ws.Range("C2:C230").Value = 45184
